# Automatische test-sync: 2025-06-23 18:35:50
#
# Adds the new "Is product Y nog op voorraad?" log entry to the Logs
# sheet, swaps the two Dashboard category labels (Offerte/Prijsaanvraag
# and Factuur/Administratie switch places) and appends the new
# "Productinformatie" tally row to the Dashboard sheet.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

# --- Logs: append the new incoming mail as row 16 -------------------------
$logs.Cells.Item(16, 1).Value = "Is product Y nog op voorraad?"
$logs.Cells.Item(16, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item(16, 3).Value = "Ik wil graag weten of product Y beschikbaar is."
$logs.Cells.Item(16, 4).Value = "Productinformatie"
$logs.Cells.Item(16, 5).Value = "Geachte klant,`nDank u wel voor uw interesse in product Y. Om u nauwkeurig te kunnen informeren over de beschikbaarheid van dit product, hebben we meer informatie nodig. Kunt u ons mogelijk het specifieke productnummer of de productnaam geven? Op die manier kunnen we direct voor u nagaan of product Y momenteel op voorraad is.`nMet vriendelijke groet,`n[Naam van het bedrijf] E-mailassistent"
$logs.Cells.Item(16, 6).Value = "2025-06-23 18:35:16"
$logs.Cells.Item(16, 7).Value = "Ja"

# --- Dashboard: the two category labels swapped places --------------------
$dash.Cells.Item(5, 1).Value = "Offerte / Prijsaanvraag"
$dash.Cells.Item(6, 1).Value = "Factuur / Administratie"

# --- Dashboard: new tally row for the Productinformatie category ----------
$dash.Cells.Item(8, 1).Value = "Productinformatie"
$dash.Cells.Item(8, 2).Value = 1

# --- Logs: extend the conditional-formatting ranges to cover row 16 -------
$dFmt = $logs.Range("D2:D15").FormatConditions.Item(1)
$dFmt.ModifyAppliesToRange($logs.Range("D2:D16"))

$gFmt = $logs.Range("G2:G15").FormatConditions.Item(1)
$gFmt.ModifyAppliesToRange($logs.Range("G2:G16"))

# --- Dashboard: point the bar chart's category/value series at row 8 too --
$chart = $dash.ChartObjects().Item(1).Chart
$series = $chart.SeriesCollection().Item(1)
$series.Formula = '=SERIES(Dashboard!$B$1,Dashboard!$A$2:$A$8,Dashboard!$B$2:$B$8,1)'
